# Tienda.xlsx update:
#   - Inventario: adjust stock counts, add new "Camisa" row
#   - Precios: add new "Camisa" pricing row
#   - Ventas: record a new day's worth of sales (2025-10-12)

$wb = $excel.ActiveWorkbook

# --- Inventario sheet: update quantities and append new product row ---
$inv = $wb.Worksheets.Item("Inventario")

$inv.Range("B2").Value = 20
$inv.Range("B4").Value = 17
$inv.Range("B5").Value = 500
$inv.Range("B6").Value = 21

$inv.Range("A7").Value = "Camisa"
$inv.Range("B7").Value = 36
$inv.Range("C7").Value = "unidades"

# --- Precios sheet: append matching pricing row for the new product ---
$pre = $wb.Worksheets.Item("Precios")

$pre.Range("A7").Value = "Camisa"
$pre.Range("B7").Value = 35000
$pre.Range("C7").Value = 72000
$pre.Range("D7").Value = 37000

# --- Ventas sheet: append new sales records ---
$ven = $wb.Worksheets.Item("Ventas")

$newSales = @(
    @(45942, "Loción",         30,  "gramos",   650,  19500),
    @(45942, "Chocolatina",    2,   "unidades", 2500, 5000),
    @(45942, "Splash",         1,   "unidades", 17000, 17000),
    @(45942, "Bolsa de Regalo",2,   "unidades", 2500, 5000),
    @(45942, "Splash",         2,   "unidades", 17000, 34000),
    @(45942, "Chocolatina",    3,   "unidades", 2500, 7500),
    @(45942, "Loción",         100, "gramos",   650,  65000)
)

$row = 22
foreach ($sale in $newSales) {
    $ven.Range("A$row").Value = $sale[0]
    $ven.Range("B$row").Value = $sale[1]
    $ven.Range("C$row").Value = $sale[2]
    $ven.Range("D$row").Value = $sale[3]
    $ven.Range("E$row").Value = $sale[4]
    $ven.Range("F$row").Value = $sale[5]

    # Carry forward the same per-column formatting used by the existing
    # sales rows (e.g. the date number format on column A) instead of
    # leaving the new rows on the sheet's default format.
    foreach ($col in @("A","B","C","D","E","F")) {
        $ven.Range("$col" + "21").Copy()
        $ven.Range("$col$row").PasteSpecial(-4122)
    }

    $row++
}
